$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("2025")
$ws1.Range("A2").Value = 53903.87531524985
$ws1.Range("B2").Value = 66631.07410080002
$ws1.Range("E2").Value = 148163.1105402212
$ws1.Range("I2").Value = 368467.1726986
$ws1.Range("M2").Value = 117062.03357365
$ws1.Range("N2").Value = 39101.99505769319
$ws1.Range("O2").Value = 69179.85095077044

$ws2 = $wb.Worksheets.Item("2030")
$ws2.Range("B2").Value = 16301.894259462
$ws2.Range("E2").Value = 252498.6354236602
$ws2.Range("I2").Value = 303793.0247600056
$ws2.Range("M2").Value = 107173.95226791
$ws2.Range("N2").Value = 62192.74864151615
$ws2.Range("O2").Value = 51778.82846662694

$ws3 = $wb.Worksheets.Item("2035")
$ws3.Range("E2").Value = 203308.5883443898
$ws3.Range("G2").Value = 36325.45083361523
$ws3.Range("I2").Value = 202590.0894332904
$ws3.Range("L2").Value = 0
$ws3.Range("M2").Value = 61279.5349989899
$ws3.Range("N2").Value = 27802.84438718831
$ws3.Range("O2").Value = 32751.34799329561
